$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (cell contents only; shared-string table is
#     rebuilt automatically by the engine) ---
$ws.Range("A1").Value = "N Terminal boundary"
$ws.Range("B1").Value = "C Terminus boundary "
$ws.Range("C1").Value = "concentration (ug/ml) "
$ws.Range("D1").Value = "yield per ml culture (ug/ml) "

# --- Re-create column D formulas as shared-formula groups, exactly as
#     Excel would do after a fill-down across each contiguous block.
#     Assigning a single relative formula string to a multi-cell Range
#     at once makes the engine emit <f t="shared" .../> groups. ---

$ws.Range("D2:D11").Formula = "=C2*80/900"
$ws.Range("D2:D11").Style = "Normal"

$ws.Range("D14:D45").Formula = "=C14*80/900"
$ws.Range("D14:D45").Style = "Normal"

$ws.Range("D46:D77").Formula = "=C46*80/900"
$ws.Range("D46:D77").Style = "Normal"

# This last shared-formula block was originally filled further down
# (through row 109) than the data currently extends (row 84), so the
# group's declared ref still spans D78:D109 while the extra trailing
# cells (D83:D109) are cleared back out again, matching the source file.
$ws.Range("D78:D109").Formula = "=C78*80/900"
$ws.Range("D78:D109").Style = "Normal"
$ws.Range("D83:D109").ClearContents()

# --- Sheet view: scroll back to the top and select A2 ---
$ws.Range("A2").Select()
